$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2507682433729189
$ws.Range("C2").Value = 139.9984142924607
$ws.Range("D2").Value = 26.60770010622114
$ws.Range("B3").Value = 0.2381696006674072
$ws.Range("C3").Value = 150.0154066158082
$ws.Range("D3").Value = 30.57983537447953
$ws.Range("B4").Value = 0.5591187777121429
$ws.Range("C4").Value = 169.9935135404647
$ws.Range("D4").Value = 45.23821949040304
$ws.Range("B5").Value = 0.2162014613825038
$ws.Range("C5").Value = 190.010250972039
$ws.Range("D5").Value = 59.8842385324773
$ws.Range("B6").Value = 0.4345211289839501
$ws.Range("C6").Value = 199.9506357312878
$ws.Range("D6").Value = 101.0966366979837
$ws.Range("B7").Value = 0.362585117781934
$ws.Range("C7").Value = 219.9635952686818
$ws.Range("D7").Value = 119.692488131674
$ws.Range("B8").Value = 0.6720973592288922
$ws.Range("C8").Value = 229.9760519082197
$ws.Range("D8").Value = 159.6033326279629
$ws.Range("B9").Value = 1.120178800818433
$ws.Range("C9").Value = 239.9592720341348
$ws.Range("D9").Value = 199.5140432874598
$ws.Range("E9").Value = 10.02459394719106
$ws.Range("F9").Value = 107.428872800964
$ws.Range("B10").Value = 0.9937320674031552
$ws.Range("C10").Value = 260.0090442692743
$ws.Range("D10").Value = 239.4739650880558
$ws.Range("E10").Value = 11.01128368342348
$ws.Range("F10").Value = 39.72116572418403
$ws.Range("B11").Value = 1.055731506605515
$ws.Range("C11").Value = 279.9618428123343
$ws.Range("D11").Value = 266.0068096765842
$ws.Range("E11").Value = 11.93321570346827
$ws.Range("F11").Value = 83.17479577808436
$ws.Range("B12").Value = 0.7433042895197757
$ws.Range("C12").Value = 290.029155460542
$ws.Range("D12").Value = 279.3085890012551
$ws.Range("E12").Value = 12.98284862508489
$ws.Range("F12").Value = 128.0326565291527
$ws.Range("B13").Value = 0.8897395155988221
$ws.Range("C13").Value = 292.02530816954
$ws.Range("D13").Value = 266.0875395603134
$ws.Range("E13").Value = 13.98355206655079
$ws.Range("F13").Value = 156.8073350724999
$ws.Range("B14").Value = 0.700772150331292
$ws.Range("C14").Value = 279.9572568635565
$ws.Range("D14").Value = 248.7352580026164
$ws.Range("E14").Value = 14.98496965002514
$ws.Range("F14").Value = 145.7498599921738
$ws.Range("B15").Value = 0.7840141341245218
$ws.Range("C15").Value = 259.929675378562
$ws.Range("D15").Value = 234.1179635629287
$ws.Range("E15").Value = 15.99151689977768
$ws.Range("F15").Value = 170.0047751077827
$ws.Range("B16").Value = 0.9788096679358043
$ws.Range("C16").Value = 252.3987584964777
$ws.Range("D16").Value = 219.4662726748273
$ws.Range("E16").Value = 16.979276812169
$ws.Range("F16").Value = 165.7068221548472
$ws.Range("B17").Value = 1.058961198360017
$ws.Range("C17").Value = 251.1828820852019
$ws.Range("D17").Value = 219.4818665510647
$ws.Range("E17").Value = 17.9896355840418
$ws.Range("F17").Value = 130.2515972013716
$ws.Range("B18").Value = 1.073699590780893
$ws.Range("C18").Value = 247.9861182083281
$ws.Range("D18").Value = 218.1831791621452
$ws.Range("E18").Value = 19.01136481086589
$ws.Range("B19").Value = 1.149208626085906
$ws.Range("C19").Value = 246.334776456127
$ws.Range("D19").Value = 214.1341865826534
$ws.Range("E19").Value = 19.99436379437963
$ws.Range("B20").Value = 1.221161738781391
$ws.Range("C20").Value = 245.61328035039
$ws.Range("D20").Value = 212.7907394884319
$ws.Range("E20").Value = 21.05292847795123
$ws.Range("B21").Value = 1.06616396739849
$ws.Range("C21").Value = 244.0086602325569
$ws.Range("D21").Value = 199.5298871748793
$ws.Range("B22").Value = 0.9490997908407564
$ws.Range("C22").Value = 240.0228911109342
$ws.Range("D22").Value = 172.9531448331074
$ws.Range("B23").Value = 0.9401319404023937
$ws.Range("C23").Value = 220.0225668827591
$ws.Range("D23").Value = 106.4419472085371
$ws.Range("B24").Value = 0.3043630684718461
$ws.Range("C24").Value = 179.9890088783918
$ws.Range("D24").Value = 66.50901503715228
$ws.Range("B25").Value = 0.3976597549510634
$ws.Range("C25").Value = 159.9767571910015
$ws.Range("D25").Value = 39.92697972083772